$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.559.16"
$ws.Range("E2").Value = "  +2.64%  "

$ws.Range("D3").Value = "2.315.11"
$ws.Range("E3").Value = "  +1.73%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'311.41"
$ws.Range("E5").Value = "  +1.53%  "

$ws.Range("D6").Value = "'102.69"
$ws.Range("E6").Value = "  +5.06%  "

$ws.Range("E7").Value = "  +1.69%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("E9").Value = "  +7.41%  "

$ws.Range("D10").Value = "'35.90"
$ws.Range("E10").Value = "  +1.22%  "

$ws.Range("E11").Value = "  +3.17%  "

$ws.Range("E12").Value = "  -0.44%  "

$ws.Range("D13").Value = "'7.04"
$ws.Range("E13").Value = "  +1.90%  "

$ws.Range("D14").Value = "2.673.80"
$ws.Range("E14").Value = "  +1.72%  "

$ws.Range("D15").Value = "'15.06"
$ws.Range("E15").Value = "  +1.37%  "

$ws.Range("D16").Value = "2.316.47"
$ws.Range("E16").Value = "  +1.70%  "

$ws.Range("D17").Value = "'0.813"
$ws.Range("E17").Value = "  +1.92%  "

$ws.Range("D18").Value = "43.453.60"
$ws.Range("E18").Value = "  +2.78%  "

$ws.Range("E19").Value = "  -0.23%  "

$ws.Range("D20").Value = "0.0₃0931"
$ws.Range("E20").Value = "  +2.43%  "

$ws.Range("E21").Value = "  +1.74%  "

$ws.Range("D22").Value = "'68.30"
$ws.Range("E22").Value = "  -0.01%  "

$ws.Range("D23").Value = "'242.11"
$ws.Range("E23").Value = "  +1.29%  "

$ws.Range("D24").Value = "'2.04"
$ws.Range("E24").Value = "  +4.37%  "

$ws.Range("E25").Value = "  +1.22%  "

$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("D28").Value = "'24.86"
$ws.Range("E28").Value = "  +4.86%  "

$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Value = "'36.90"
$ws.Range("E29").Value = "  -3.16%  "

$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Value = "'9.65"
$ws.Range("E30").Value = "  +1.13%  "

$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "'2.12"
$ws.Range("E31").Value = "  -0.11%  "

$ws.Range("D32").Value = "'167.89"
$ws.Range("E32").Value = "  +3.77%  "

$ws.Range("D33").Value = "'5.30"
$ws.Range("E33").Value = "  +0.74%  "

$ws.Range("E34").Value = "  +0.08%  "

$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "'2.53"
$ws.Range("E35").Value = "  +6.77%  "

$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").Value = "'3.12"
$ws.Range("E36").Value = "  -2.47%  "

$ws.Range("D37").Value = "'0.0746"
$ws.Range("E37").Value = "  +0.78%  "

$ws.Range("E38").Value = "  +0.02%  "

$ws.Range("E39").Value = "  +1.88%  "

$ws.Range("E40").Value = "  +2.61%  "

$ws.Range("E41").Value = "  +1.75%  "

$ws.Range("D42").Value = "'4.33"
$ws.Range("E42").Value = "  +5.23%  "

$ws.Range("E43").Value = "  -0.55%  "

$ws.Range("D44").Value = "'19.41"
$ws.Range("E44").Value = "  +0.85%  "

$ws.Range("E45").Value = "  +2.71%  "

$ws.Range("D46").Value = "1.972.55"
$ws.Range("E46").Value = "  +1.03%  "

$ws.Range("D47").Value = "'2.99"
$ws.Range("E47").Value = "  +2.87%  "

$ws.Range("D48").Value = "'9.94"
$ws.Range("E48").Value = "  -0.38%  "

$ws.Range("D49").Value = "'55.64"
$ws.Range("E49").Value = "  +3.47%  "

$ws.Range("D50").Value = "'2.93"
$ws.Range("E50").Value = "  +5.14%  "

$ws.Range("D51").Value = "'1.58"
$ws.Range("E51").Value = "  +6.82%  "
